$d = $word.ActiveDocument

function Get-ParagraphIndexContaining($doc, $text) {
    $idx = 0
    foreach ($p in $doc.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text.Contains($text)) {
            return $idx
        }
    }
    return -1
}

# --- Step 1: merge the two runs of "Caso concreto..." into a single run. ---
$d.Content.Find.Execute(
    "Caso concreto: film vietati a minori di 18 anni e controllo dell’età con verifiche.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Caso concreto: film vietati a minori di 18 anni e controllo dell’età con verifiche.",
    2)

# --- Step 2: remove the "Per ChatGPT" heading block (4 paragraphs):
#     "Per ChatGPT", "Questo il codice del contratto...", and two blank paragraphs.
$iStart = Get-ParagraphIndexContaining $d "Per ChatGPT"
$iEnd = Get-ParagraphIndexContaining $d "Questo il codice del contratto"
$iEnd = $iEnd + 2
$pStart = $d.Paragraphs($iStart)
$pEnd = $d.Paragraphs($iEnd)
$r = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$r.Delete()

# --- Step 3: remove the page-break paragraph and the subsequent block of
#     paragraphs describing DID generation, through the end of the document.
$iStart2 = Get-ParagraphIndexContaining $d "La generazione del DID dovrebbe avvenire"
$count = $d.Paragraphs.Count
$pStart2 = $d.Paragraphs($iStart2)
$pEnd2 = $d.Paragraphs($count)
$r2 = $d.Range($pStart2.Range.Start, $pEnd2.Range.End)
$r2.Delete()
